# Update the results for the classifiers:
# - Populate the previously-empty "PCA-300-Corpus" sheet with the
#   PCA-300 corpus-level classifier results (headers, ngram labels and
#   accuracy figures), copying the layout/formatting from the sibling
#   "PCA-300-Polarity" sheet.
# - Make "PCA-300-Corpus" the active/selected sheet, replacing
#   "PCA-300-Polarity" as the tab that was previously selected.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("PCA-300-Polarity")
$dst = $wb.Worksheets.Item("PCA-300-Corpus")

# Bring over the header row / ngram labels and cell formatting (borders,
# wrap text, row heights, etc.) from the already-filled-in Polarity sheet.
$src.Range("A1:M4").Copy($dst.Range("A1"))

# The source sheet bolds the winning train/test columns (J2:K2) with a
# different style; the corpus sheet does not, so normalize those two
# cells back to the plain data style used by the rest of the row (copy
# formatting only, values are overwritten below anyway).
$dst.Range("L2").Copy()
$dst.Range("J2:K2").PasteSpecial(-4122)

# Write in the updated classifier results for the corpus run.
$dst.Range("B2").Value = 0.69544600000000001
$dst.Range("C2").Value = 0.69476899999999997
$dst.Range("D2").Value = 0.54695400000000005
$dst.Range("E2").Value = 0.549597
$dst.Range("F2").Value = 0.69201400000000002
$dst.Range("G2").Value = 0.69125599999999998
$dst.Range("H2").Value = 0.55613699999999999
$dst.Range("I2").Value = 0.54657100000000003
$dst.Range("J2").Value = 0.71076099999999998
$dst.Range("K2").Value = 0.68569100000000005
$dst.Range("L2").Value = 0.59662199999999999
$dst.Range("M2").Value = 0.59303700000000004

$dst.Range("B3").Value = 0.55484999999999995
$dst.Range("C3").Value = 0.55728299999999997
$dst.Range("D3").Value = 0.51846899999999996
$dst.Range("E3").Value = 0.52257200000000004
$dst.Range("F3").Value = 0.55385300000000004
$dst.Range("G3").Value = 0.55655299999999996
$dst.Range("H3").Value = 0.53728500000000001
$dst.Range("I3").Value = 0.53923200000000004
$dst.Range("J3").Value = 0.589005
$dst.Range("K3").Value = 0.56747400000000003
$dst.Range("L3").Value = 0.53417800000000004
$dst.Range("M3").Value = 0.53787600000000002

$dst.Range("B4").Value = 0.48330499999999998
$dst.Range("C4").Value = 0.47436699999999998
$dst.Range("D4").Value = 0.47007700000000002
$dst.Range("E4").Value = 0.46372400000000003
$dst.Range("F4").Value = 0.48258600000000001
$dst.Range("G4").Value = 0.47429700000000002
$dst.Range("H4").Value = 0.48403600000000002
$dst.Range("I4").Value = 0.47067999999999999
$dst.Range("J4").Value = 0.49601200000000001
$dst.Range("K4").Value = 0.477219
$dst.Range("L4").Value = 0.47928199999999999
$dst.Range("M4").Value = 0.470611

# Match the (taller, wrapped-header) row heights used on the source sheet.
$dst.Rows.Item(1).RowHeight = 46
$dst.Rows.Item(2).RowHeight = 17
$dst.Rows.Item(3).RowHeight = 17
$dst.Rows.Item(4).RowHeight = 17

# Move the selection on the Polarity sheet first (selecting there makes
# it active momentarily), then finish by activating/selecting on the
# Corpus sheet so it ends up as the selected tab.
$src.Range("K9").Select()

$dst.Activate()
$dst.Range("A4").Select()
